# Updates crypto Price (D) and Volume(1h) (E) columns to refreshed
# values, matching a new data pull from the source feed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '36.760.14'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +2.02%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.963.19'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +2.53%  '

$ws.Range("E4").Value = '  -0.15%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '244.00'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.51%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.615'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +2.25%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '58.49'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +5.75%  '

$ws.Range("E8").Value = '  -0.05%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0811'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -0.44%  '

$ws.Range("E11").Value = '  +0.45%  '

$ws.Range("E12").Value = '  +6.88%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.246.02'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +2.24%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.821'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +2.46%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '13.71'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +4.59%  '

$ws.Range("E16").Value = '  +2.72%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.960.76'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +1.15%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '36.713.46'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +2.01%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '69.74'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +2.14%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0863'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +1.92%  '

$ws.Range("E21").Value = '  +4.02%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '228.17'
$ws.Range("D22").ClearFormats()

$ws.Range("E23").Value = '  -0.17%  '

$ws.Range("E24").Value = '  -0.80%  '

$ws.Range("E25").Value = '  +3.38%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.30'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +1.18%  '

$ws.Range("E27").Value = '  +16.21%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '161.00'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -0.24%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.37'
$ws.Range("D29").ClearFormats()

$ws.Range("E30").Value = '  +2.24%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.12'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.38%  '

$ws.Range("E32").Value = '  +2.01%  '

$ws.Range("E33").Value = '  +0.19%  '

$ws.Range("E34").Value = '  +0.22%  '

$ws.Range("E35").Value = '  +6.56%  '

$ws.Range("E36").Value = '  -0.17%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.44'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +21.70%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.20'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +4.52%  '

$ws.Range("E39").Value = '  -0.86%  '

$ws.Range("E40").Value = '  +5.59%  '

$ws.Range("E41").Value = '  +2.39%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0213'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +3.82%  '

$ws.Range("E43").Value = '  +1.37%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.04'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +4.40%  '

$ws.Range("E45").Value = '  +2.59%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.346.24'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +1.36%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '87.45'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +1.27%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.12'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.21%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.84'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +1.74%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.138.05'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +2.30%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '43.55'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -2.32%  '
